$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing forecast-error rows (data for Q0..Q8, currently in rows 2-10)
# down by one row (into rows 3-11), as a new observation becomes available
# and the naive component forecaster error table picks up one more data point.
# Process from the bottom row upward so we don't overwrite values before they
# are copied down. Read via Value2 (reliable for reads in this runtime) and
# write via Value.
for ($r = 10; $r -ge 2; $r--) {
    $b = $ws.Cells.Item($r, 2).Value2
    $c = $ws.Cells.Item($r, 3).Value2
    $d = $ws.Cells.Item($r, 4).Value2
    $e = $ws.Cells.Item($r, 5).Value2
    $f = $ws.Cells.Item($r, 6).Value2

    $ws.Cells.Item($r + 1, 2).Value = $b
    $ws.Cells.Item($r + 1, 3).Value = $c
    $ws.Cells.Item($r + 1, 4).Value = $d
    $ws.Cells.Item($r + 1, 5).Value = $e
    $ws.Cells.Item($r + 1, 6).Value = $f
}

# Row 2 (Q0) now gets freshly computed error-statistic values.
$ws.Cells.Item(2, 2).Value = -0.00515214134517856
$ws.Cells.Item(2, 3).Value = 0.6651924495933969
$ws.Cells.Item(2, 4).Value = 0.7443617204222924
$ws.Cells.Item(2, 5).Value = 0.8627640004209103
$ws.Cells.Item(2, 6).Value = 0.8863899424499349

# Column G holds the sample-size N used for each row; it grows by one for
# every row because a new data point has been matched/added.
for ($r = 2; $r -le 11; $r++) {
    $n = $ws.Cells.Item($r, 7).Value2
    $ws.Cells.Item($r, 7).Value = $n + 1
}
